$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected FilesTab ICDC query: drop the "File Type" and "Breed" return
# columns (matches the "corrected ICDC Breed 1-14 scripts" commit).
$newFilesQuery = @"
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
WHERE demo.breed IN ['Cocker Spaniel']
OPTIONAL MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
WITH DISTINCT f, parent, c, demo, diag, s
RETURN  coalesce(f.file_name, '') AS ``File Name``,
         coalesce(labels(parent)[0], '') AS ``Association``,
        coalesce(f.file_description, '') AS ``Description``,
        coalesce(f.file_format, '') AS ``Format``,
        coalesce(f.file_size, '') AS ``Size``,
        coalesce(c.case_id, '') AS ``Case ID``,
         coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS ``Study Code``
"@

$ws.Range("B4").Value = $newFilesQuery

# Row 4 shrank (two fewer wrapped lines) once the query text got shorter.
$ws.Rows(4).RowHeight = 217.5

# The author's selection ended up on the cell they just edited, scrolled
# so row 4 is visible at the top of the pane.
$null = $ws.Range("B4").Select()
$excel.ActiveWindow.ScrollRow = 4
$excel.ActiveWindow.ScrollColumn = 1
